# events API TODO.docx edit script
# Implements:
#  - Title text change: "events API TODO" -> "Events API (events-api) TODO"
#    with the _GoBack bookmark relocated to sit right after "(events-api"
#    and the trailing line break removed from the title paragraph.
#  - Three new paragraphs inserted right after the title:
#      1) an empty paragraph
#      2) a paragraph containing a hyperlink to
#         https://kb.novaordis.com/index.php/Events-api followed by a space
#      3) an empty paragraph
#  - The existing "Keep ... up to date." paragraph: the old _GoBack bookmark
#    is removed from it (it moved to the title) and its two-run hyperlink
#    ("https://kb.novaordis.com/index.php" + "/" + "Events-api_Concepts")
#    is consolidated into a single-run hyperlink with the full URL text.
#  - One additional empty paragraph appended at the very end of the document.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Title paragraph: change text and drop the trailing line break.
# ---------------------------------------------------------------------
$d.Content.Find.Execute("events API TODO", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Events API (events-api) TODO", 2)

$titlePara = $d.Paragraphs(1)
$titleRange = $titlePara.Range
# The paragraph now ends with "...TODO" + <w:br/> (vertical-tab, char 11) +
# the paragraph mark (char 13). Remove just the line-break character.
$breakRange = $d.Range($titleRange.End - 2, $titleRange.End - 1)
$breakRange.Delete()

# Move the (single, document-wide) _GoBack bookmark so that it sits right
# after "(events-api" in the title text.
$findRange = $d.Content
$findRange.Find.Execute("(events-api", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "", 0)
$bookmarkPos = $findRange.End
$bmRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# ---------------------------------------------------------------------
# 2. Insert three new paragraphs right after the title paragraph, before
#    the "JBoss-Related Heuristics" list item.
# ---------------------------------------------------------------------
$jbossPara = $d.Paragraphs(2)
$insertHere = $jbossPara.Range
$insertHere.Collapse(1)

# paragraph A (empty)
$insertHere.InsertParagraphBefore()
# paragraph B (hyperlink)
$insertHere.InsertParagraphBefore()
# paragraph C (empty)
$insertHere.InsertParagraphBefore()

# After three inserts, paragraphs 2, 3 and 4 are the new ones (in order),
# and "JBoss-Related Heuristics" becomes paragraph 5.
$paraA = $d.Paragraphs(2)
$paraB = $d.Paragraphs(3)
$paraC = $d.Paragraphs(4)

foreach ($p in @($paraA, $paraB, $paraC)) {
    $p.Range.ListFormat.RemoveNumbers()
    $p.Style = $d.Styles("Normal")
    $p.Range.Font.Name = "Garamond"
}

# Fill paragraph B with the hyperlink "https://kb.novaordis.com/index.php/Events-api"
# followed by a trailing space.
$bStart = $paraB.Range.Start
$bEnd = $paraB.Range.End - 1
$bTextRange = $d.Range($bStart, $bEnd)
$bTextRange.Text = " "

$linkUrl = "https://kb.novaordis.com/index.php/Events-api"
$linkInsertPoint = $d.Range($bStart, $bStart)
$newHyperlink = $d.Hyperlinks.Add($linkInsertPoint, $linkUrl, "", "", $linkUrl)
$newHyperlink.Range.Font.Name = "Garamond"

# ---------------------------------------------------------------------
# 3. Fix up the "Keep ... up to date." paragraph: consolidate the
#    hyperlink runs into a single run with the full URL text.
# ---------------------------------------------------------------------
$keepHyperlink = $d.Hyperlinks(1)
$hStart = $keepHyperlink.Range.Start
$hEnd = $keepHyperlink.Range.End
$hDeleteRange = $d.Range($hStart, $hEnd)
$hDeleteRange.Delete()

$concUrl = "https://kb.novaordis.com/index.php/Events-api_Concepts"
$concInsertPoint = $d.Range($hStart, $hStart)
$concInsertPoint.InsertAfter($concUrl)
$concRange = $d.Range($hStart, $hStart + $concUrl.Length)
$newHyperlink2 = $d.Hyperlinks.Add($concRange, $concUrl, "", "", $concUrl)
$newHyperlink2.Range.Font.Name = "Garamond"

# ---------------------------------------------------------------------
# 4. Append one more empty paragraph at the very end of the document.
# ---------------------------------------------------------------------
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($lastIndex)
$endRange = $lastPara.Range
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
